$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = -21.9213
$ws.Range("A13").Value = -22.14900000000001
$ws.Range("A16").Value = -21.59359999999999
$ws.Range("A18").Value = -22.22370000000001
$ws.Range("A20").Value = -21.03339999999998
$ws.Range("A26").Value = -21.01259999999997
$ws.Range("A27").Value = -21.68009999999999
$ws.Range("A29").Value = -20.79129999999997
$ws.Range("A35").Value = -20.99569999999999
$ws.Range("A36").Value = -21.0463
$ws.Range("A45").Value = -21.52959999999999
$ws.Range("A55").Value = -22.3012
$ws.Range("A57").Value = -22.25360000000001
$ws.Range("A69").Value = -21.67099999999999
$ws.Range("A76").Value = -19.61219999999998
$ws.Range("A78").Value = -19.67819999999998
$ws.Range("A82").Value = -22.17670000000001
$ws.Range("A83").Value = -21.98520000000001
$ws.Range("A93").Value = -20.72019999999997
$ws.Range("A97").Value = -22.00350000000001
